$wb = $excel.ActiveWorkbook

# Sheet ALC, row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3413.4
$ws.Range("I19").Value = 6750.5
$ws.Range("J19").Value = 1188.6666
$ws.Range("K19").Value = 6750.5
$ws.Range("L19").Value = 1188.6666
$ws.Range("M19").Value = -6575.5
$ws.Range("N19").Value = -1538.6666

# Sheet ALC, row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1049.0869
$ws.Range("I28").Value = 929.375
$ws.Range("K28").Value = 929.375
$ws.Range("M28").Value = -444.375

# Sheet ALC, row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 997.7742
$ws.Range("I107").Value = 1139.1666
$ws.Range("K107").Value = 1139.1666
$ws.Range("M107").Value = 780.8334

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1797.6552
$ws.Range("I137").Value = 1055.091
$ws.Range("J137").Value = 4131.4287
$ws.Range("K137").Value = 3165.273
$ws.Range("L137").Value = 12394.2861
$ws.Range("M137").Value = -615.2729999999997
$ws.Range("N137").Value = -17494.2861

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2385.42
$ws.Range("J138").Value = 2732.077
$ws.Range("L138").Value = 8196.231
$ws.Range("N138").Value = -18476.231

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 19500
$ws.Range("I141").Value = 35000
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 105000
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -99820
$ws.Range("N141").Value = -22360

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 53483.58
$ws.Range("I45").Value = 84050.086
$ws.Range("K45").Value = 84050.086
$ws.Range("M45").Value = -83673.086

# Sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 862.5185
$ws.Range("I110").Value = 775.04346
$ws.Range("J110").Value = 1365.5
$ws.Range("K110").Value = 775.04346
$ws.Range("L110").Value = 1365.5
$ws.Range("M110").Value = 1269.95654
$ws.Range("N110").Value = -5455.5

# Sheet ARM, row 117
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 32400
$ws.Range("J117").Value = 32400
$ws.Range("L117").Value = 32400
$ws.Range("N117").Value = -41578

# Sheet ARM, row 118
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H118").Value = 30971.428
$ws.Range("J118").Value = 30971.428
$ws.Range("L118").Value = 30971.428
$ws.Range("N118").Value = -34285.428

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1373.6666
$ws.Range("I122").Value = 1437.5
$ws.Range("J122").Value = 1150.25
$ws.Range("K122").Value = 4312.5
$ws.Range("L122").Value = 3450.75
$ws.Range("M122").Value = -1862.5
$ws.Range("N122").Value = -8350.75

# Sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1982.7778
$ws.Range("I105").Value = 2031.25
$ws.Range("J105").Value = 1595
$ws.Range("K105").Value = 2031.25
$ws.Range("L105").Value = 1595
$ws.Range("M105").Value = -284.25
$ws.Range("N105").Value = -5089

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1267.8889
$ws.Range("I107").Value = 1242.4286
$ws.Range("J107").Value = 1357
$ws.Range("K107").Value = 1242.4286
$ws.Range("L107").Value = 1357
$ws.Range("M107").Value = 677.5714
$ws.Range("N107").Value = -5197

# Sheet BSM, row 137
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 107500
$ws.Range("J137").Value = 107500
$ws.Range("L137").Value = 107500
$ws.Range("N137").Value = -117700

# Sheet BSM, row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 45499.832
$ws.Range("J140").Value = 45499.832
$ws.Range("L140").Value = 45499.832
$ws.Range("N140").Value = -55859.832

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1241.4286
$ws.Range("I16").Value = 1290.2222
$ws.Range("J16").Value = 1153.6
$ws.Range("K16").Value = 1290.2222
$ws.Range("L16").Value = 1153.6
$ws.Range("M16").Value = -1003.2222
$ws.Range("N16").Value = -1727.6

# Sheet CRP, row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1009.5
$ws.Range("I105").Value = 1027.7142
$ws.Range("J105").Value = 967
$ws.Range("K105").Value = 1027.7142
$ws.Range("L105").Value = 967
$ws.Range("M105").Value = 719.2858000000001
$ws.Range("N105").Value = -4461

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 424.04544
$ws.Range("I107").Value = 404.66666
$ws.Range("J107").Value = 511.25
$ws.Range("K107").Value = 404.66666
$ws.Range("L107").Value = 511.25
$ws.Range("M107").Value = 1515.33334
$ws.Range("N107").Value = -4351.25

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1241.4286
$ws.Range("I113").Value = 1290.2222
$ws.Range("J113").Value = 1153.6
$ws.Range("K113").Value = 1290.2222
$ws.Range("L113").Value = 1153.6
$ws.Range("M113").Value = 879.7778000000001
$ws.Range("N113").Value = -5493.6

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1211.0869
$ws.Range("I5").Value = 1303.125
$ws.Range("J5").Value = 1000.7143
$ws.Range("K5").Value = 3909.375
$ws.Range("L5").Value = 3002.1429
$ws.Range("M5").Value = -3797.375
$ws.Range("N5").Value = -3226.1429

# Sheet CUL, row 31
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 2900
$ws.Range("J31").Value = 2900
$ws.Range("L31").Value = 8700
$ws.Range("N31").Value = -9276

# Sheet CUL, row 76
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 1850
$ws.Range("I76").Value = 1166.6666
$ws.Range("J76").Value = 3900
$ws.Range("K76").Value = 3499.9998
$ws.Range("L76").Value = 11700
$ws.Range("M76").Value = -3116.9998
$ws.Range("N76").Value = -12466

# Sheet CUL, row 79
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 1850
$ws.Range("I79").Value = 1166.6666
$ws.Range("J79").Value = 3900
$ws.Range("K79").Value = 3499.9998
$ws.Range("L79").Value = 11700
$ws.Range("M79").Value = -2173.9998
$ws.Range("N79").Value = -14352

# Sheet CUL, row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 1919.5
$ws.Range("J99").Value = 3600
$ws.Range("L99").Value = 10800
$ws.Range("N99").Value = -15292

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1335601.1
$ws.Range("I131").Value = 7001.1113
$ws.Range("J131").Value = 1516773.9
$ws.Range("K131").Value = 21003.3339
$ws.Range("L131").Value = 4550321.699999999
$ws.Range("M131").Value = -15963.3339
$ws.Range("N131").Value = -4560401.699999999

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1211.0869
$ws.Range("I135").Value = 1303.125
$ws.Range("J135").Value = 1000.7143
$ws.Range("K135").Value = 11728.125
$ws.Range("L135").Value = 9006.4287
$ws.Range("M135").Value = -9193.125
$ws.Range("N135").Value = -14076.4287

# Sheet GSM, row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1376.7297
$ws.Range("I97").Value = 1441.4073
$ws.Range("J97").Value = 1202.1
$ws.Range("K97").Value = 1441.4073
$ws.Range("L97").Value = 1202.1
$ws.Range("M97").Value = -945.4073000000001
$ws.Range("N97").Value = -2194.1

# Sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2186.8462
$ws.Range("I113").Value = 2092
$ws.Range("K113").Value = 2092
$ws.Range("M113").Value = 78

# Sheet LTW, row 45
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 14448.5
$ws.Range("I45").Value = 14499
$ws.Range("J45").Value = 14398
$ws.Range("K45").Value = 14499
$ws.Range("L45").Value = 14398
$ws.Range("M45").Value = -14092
$ws.Range("N45").Value = -15212

# Sheet LTW, row 48
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 11950
$ws.Range("J48").Value = 13900
$ws.Range("L48").Value = 13900
$ws.Range("N48").Value = -15222

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 684.2143
$ws.Range("I61").Value = 567.7
$ws.Range("J61").Value = 975.5
$ws.Range("K61").Value = 567.7
$ws.Range("L61").Value = 975.5
$ws.Range("M61").Value = -365.7
$ws.Range("N61").Value = -1379.5

# Sheet LTW, row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1198.2727
$ws.Range("I100").Value = 1220.1
$ws.Range("J100").Value = 980
$ws.Range("K100").Value = 1220.1
$ws.Range("L100").Value = 980
$ws.Range("M100").Value = -679.0999999999999
$ws.Range("N100").Value = -2062

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 684.2143
$ws.Range("I113").Value = 567.7
$ws.Range("J113").Value = 975.5
$ws.Range("K113").Value = 567.7
$ws.Range("L113").Value = 975.5
$ws.Range("M113").Value = 1602.3
$ws.Range("N113").Value = -5315.5

# Sheet WVR, row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 474.04
$ws.Range("I107").Value = 326
$ws.Range("J107").Value = 610.6923
$ws.Range("K107").Value = 978
$ws.Range("L107").Value = 1832.0769
$ws.Range("M107").Value = 942
$ws.Range("N107").Value = -5672.0769

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1220.5
$ws.Range("I113").Value = 1251
$ws.Range("J113").Value = 1190
$ws.Range("K113").Value = 3753
$ws.Range("L113").Value = 3570
$ws.Range("M113").Value = -1583
$ws.Range("N113").Value = -7910
